# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new trailing columns (date / legislator_name / legislator_id)
# to the "股票" (stocks) worksheet, filling every existing data row with
# the report date, the legislator's name and numeric id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "楊麗環"
$legislatorId   = 960
$reportDate     = "2011-11-22"

# Find the last used data row in column A (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Header row (H1:J1) -----------------------------------------------
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Copy the formatting of the existing header cell (G1) onto the new
# header cells so they keep the bold / centered / bordered look.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (H2:J<lastRow>) -----------------------------------------
# Format column H as text first so the "yyyy-mm-dd" string is stored
# literally instead of being auto-converted into a date serial number.
$ws.Range("H2:H$lastRow").NumberFormat = "@"
$ws.Range("H2:H$lastRow").Value = $reportDate
$ws.Range("I2:I$lastRow").Value = $legislatorName
$ws.Range("J2:J$lastRow").Value = $legislatorId

# Normalize formatting of the new columns to match the existing data
# columns (e.g. column G) and drop the temporary text number format.
$ws.Range("G2:G$lastRow").Copy()
$ws.Range("H2:H$lastRow").PasteSpecial(-4122)
$ws.Range("I2:I$lastRow").PasteSpecial(-4122)
$ws.Range("J2:J$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
